$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 569.9286
$ws.Range("I2").Value = 547
$ws.Range("J2").Value = 592.8570999999999
$ws.Range("K2").Value = 547
$ws.Range("L2").Value = 592.8570999999999
$ws.Range("M2").Value = -434
$ws.Range("N2").Value = -818.8570999999999
$ws.Range("H40").Value = 3677.1428
$ws.Range("I40").Value = 2040
$ws.Range("J40").Value = 3950
$ws.Range("K40").Value = 2040
$ws.Range("L40").Value = 3950
$ws.Range("M40").Value = -1865
$ws.Range("N40").Value = -4300
$ws.Range("H42").Value = 420.7143
$ws.Range("I42").Value = 166.66667
$ws.Range("J42").Value = 611.25
$ws.Range("K42").Value = 500.00001
$ws.Range("L42").Value = 1833.75
$ws.Range("M42").Value = -270.00001
$ws.Range("N42").Value = -2293.75
$ws.Range("H86").Value = 2331.68
$ws.Range("I86").Value = 2132.7222
$ws.Range("J86").Value = 2843.2856
$ws.Range("K86").Value = 2132.7222
$ws.Range("L86").Value = 2843.2856
$ws.Range("M86").Value = -1009.7222
$ws.Range("H89").Value = 2331.68
$ws.Range("I89").Value = 2132.7222
$ws.Range("J89").Value = 2843.2856
$ws.Range("K89").Value = 10663.611
$ws.Range("L89").Value = 14216.428
$ws.Range("M89").Value = -5047.611000000001
$ws.Range("H106").Value = 4659.125
$ws.Range("I106").Value = 3962.1667
$ws.Range("J106").Value = 6750
$ws.Range("K106").Value = 3962.1667
$ws.Range("L106").Value = 6750
$ws.Range("M106").Value = -3331.1667
$ws.Range("H113").Value = 6490.909
$ws.Range("I113").Value = 5950
$ws.Range("J113").Value = 6800
$ws.Range("K113").Value = 5950
$ws.Range("L113").Value = 6800
$ws.Range("M113").Value = -2696
$ws.Range("N113").Value = -13308
$ws.Range("H116").Value = 5188
$ws.Range("I116").Value = 5250.8335
$ws.Range("J116").Value = 4999.5
$ws.Range("K116").Value = 5250.8335
$ws.Range("L116").Value = 4999.5
$ws.Range("M116").Value = -1808.8335

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2341.9048
$ws.Range("I2").Value = 2645
$ws.Range("J2").Value = 1849.375
$ws.Range("K2").Value = 2645
$ws.Range("L2").Value = 1849.375
$ws.Range("M2").Value = -2532
$ws.Range("N2").Value = -2075.375
$ws.Range("H45").Value = 1780
$ws.Range("I45").Value = 1795
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 1795
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -1418
$ws.Range("N45").Value = -2504
$ws.Range("H74").Value = 45508944
$ws.Range("I74").Value = 62572160
$ws.Range("J74").Value = 7031.1665
$ws.Range("K74").Value = 62572160
$ws.Range("L74").Value = 7031.1665
$ws.Range("M74").Value = -62571286
$ws.Range("H77").Value = 45508944
$ws.Range("I77").Value = 62572160
$ws.Range("J77").Value = 7031.1665
$ws.Range("K77").Value = 312860800
$ws.Range("L77").Value = 35155.8325
$ws.Range("M77").Value = -312856432
$ws.Range("H116").Value = 2341.9048
$ws.Range("I116").Value = 2645
$ws.Range("J116").Value = 1849.375
$ws.Range("K116").Value = 2645
$ws.Range("L116").Value = 1849.375
$ws.Range("M116").Value = -351
$ws.Range("N116").Value = -6437.375
$ws.Range("H132").Value = 21278316
$ws.Range("I132").Value = 1735.25
$ws.Range("J132").Value = 333334850
$ws.Range("K132").Value = 5205.75
$ws.Range("L132").Value = 1000004550
$ws.Range("M132").Value = -2675.75
$ws.Range("N132").Value = -1000009610

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2341.9048
$ws.Range("I3").Value = 2645
$ws.Range("J3").Value = 1849.375
$ws.Range("K3").Value = 2645
$ws.Range("L3").Value = 1849.375
$ws.Range("M3").Value = -2531
$ws.Range("N3").Value = -2077.375
$ws.Range("H29").Value = 4999.6665
$ws.Range("I29").Value = 4999.5
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 4999.5
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = -4710.5
$ws.Range("N29").Value = -5578
$ws.Range("H120").Value = 68000
$ws.Range("I120").Value = 68000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 68000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -63162
$ws.Range("N120").ClearContents()
$ws.Range("H134").Value = 2815.6511
$ws.Range("I134").Value = 2739.8333
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 8219.499899999999
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -5684.499899999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3843.4827
$ws.Range("I22").Value = 5086.095
$ws.Range("J22").Value = 581.625
$ws.Range("K22").Value = 5086.095
$ws.Range("L22").Value = 581.625
$ws.Range("M22").Value = -4736.095
$ws.Range("N22").Value = -1281.625
$ws.Range("H134").Value = 1309.2632
$ws.Range("I134").Value = 1215.3889
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 3646.1667
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -1111.1667
$ws.Range("H141").Value = 444466.4
$ws.Range("I141").Value = 52666.668
$ws.Range("J141").Value = 575066.3
$ws.Range("K141").Value = 52666.668
$ws.Range("L141").Value = 575066.3
$ws.Range("M141").Value = -47486.668
$ws.Range("N141").Value = -585426.3

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 66867.60000000001
$ws.Range("I9").Value = 94468
$ws.Range("J9").Value = 2466.6667
$ws.Range("K9").Value = 283404
$ws.Range("L9").Value = 7400.000100000001
$ws.Range("M9").Value = -283180
$ws.Range("H18").Value = 2279.5833
$ws.Range("I18").Value = 1009.2857
$ws.Range("J18").Value = 4058
$ws.Range("K18").Value = 3027.8571
$ws.Range("L18").Value = 12174
$ws.Range("M18").Value = -2858.8571
$ws.Range("N18").Value = -12512
$ws.Range("H62").Value = 9527533
$ws.Range("I62").Value = 1256
$ws.Range("J62").Value = 13338043
$ws.Range("K62").Value = 3768
$ws.Range("L62").Value = 40014129
$ws.Range("M62").Value = -3082
$ws.Range("N62").Value = -40015501
$ws.Range("H63").Value = 4628.6665
$ws.Range("I63").Value = 4443
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 13329
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -12580
$ws.Range("H64").Value = 22979
$ws.Range("I64").Value = 1900
$ws.Range("J64").Value = 28248.75
$ws.Range("K64").Value = 5700
$ws.Range("L64").Value = 84746.25
$ws.Range("M64").Value = -5430
$ws.Range("H65").Value = 9527533
$ws.Range("I65").Value = 1256
$ws.Range("J65").Value = 13338043
$ws.Range("K65").Value = 11304
$ws.Range("L65").Value = 120042387
$ws.Range("M65").Value = -7872
$ws.Range("N65").Value = -120049251
$ws.Range("H66").Value = 4628.6665
$ws.Range("I66").Value = 4443
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 39987
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -36243
$ws.Range("H67").Value = 22979
$ws.Range("I67").Value = 1900
$ws.Range("J67").Value = 28248.75
$ws.Range("K67").Value = 5700
$ws.Range("L67").Value = 84746.25
$ws.Range("M67").Value = -4764
$ws.Range("H68").Value = 1497.5
$ws.Range("H71").Value = 1497.5
$ws.Range("H98").Value = 2028.1428
$ws.Range("I98").Value = 955
$ws.Range("J98").Value = 2320.818
$ws.Range("K98").Value = 2865
$ws.Range("L98").Value = 6962.454000000001
$ws.Range("M98").Value = -1367
$ws.Range("N98").Value = -9958.454000000002
$ws.Range("H99").Value = 1645.2858
$ws.Range("I99").Value = 1103
$ws.Range("J99").Value = 3001
$ws.Range("K99").Value = 3309
$ws.Range("L99").Value = 9003
$ws.Range("M99").Value = -1063
$ws.Range("N99").Value = -13495
$ws.Range("H100").Value = 200
$ws.Range("I100").Value = 200
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 211
$ws.Range("N100").ClearContents()
$ws.Range("H102").Value = 1999
$ws.Range("I102").Value = 1999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5997
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3563
$ws.Range("H103").Value = 1741.6666
$ws.Range("I103").Value = 600
$ws.Range("J103").Value = 1970
$ws.Range("K103").Value = 1800
$ws.Range("L103").Value = 5910
$ws.Range("M103").Value = -921
$ws.Range("N103").Value = -7668
$ws.Range("H104").Value = 1761.6923
$ws.Range("I104").Value = 918
$ws.Range("J104").Value = 2014.8
$ws.Range("K104").Value = 2754
$ws.Range("L104").Value = 6044.4
$ws.Range("M104").Value = -133
$ws.Range("N104").Value = -11286.4
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 13785.429
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 13785.429
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 41356.287
$ws.Range("N106").Value = -43248.287
$ws.Range("M106").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 14250
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 14250
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 14250
$ws.Range("N106").Value = -16774

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3453.2778
$ws.Range("I7").Value = 3368.7646
$ws.Range("J7").Value = 4890
$ws.Range("K7").Value = 3368.7646
$ws.Range("L7").Value = 4890
$ws.Range("M7").Value = -3256.7646
$ws.Range("N7").Value = -5114
$ws.Range("H40").Value = 4508.875
$ws.Range("I40").Value = 4153.7144
$ws.Range("J40").Value = 6995
$ws.Range("K40").Value = 4153.7144
$ws.Range("L40").Value = 6995
$ws.Range("M40").Value = -4017.7144
$ws.Range("H46").Value = 2250.1082
$ws.Range("I46").Value = 1042.8096
$ws.Range("J46").Value = 3834.6875
$ws.Range("K46").Value = 1042.8096
$ws.Range("L46").Value = 3834.6875
$ws.Range("M46").Value = -854.8096
$ws.Range("N46").Value = -4210.6875
$ws.Range("H55").Value = 639.3077
$ws.Range("I55").Value = 351.94116
$ws.Range("J55").Value = 1182.1111
$ws.Range("K55").Value = 351.94116
$ws.Range("L55").Value = 1182.1111
$ws.Range("M55").Value = -178.94116
$ws.Range("H61").Value = 4488.6924
$ws.Range("I61").Value = 2240
$ws.Range("J61").Value = 5894.125
$ws.Range("K61").Value = 2240
$ws.Range("L61").Value = 5894.125
$ws.Range("M61").Value = -2038
$ws.Range("H113").Value = 4488.6924
$ws.Range("I113").Value = 2240
$ws.Range("J113").Value = 5894.125
$ws.Range("K113").Value = 2240
$ws.Range("L113").Value = 5894.125
$ws.Range("M113").Value = -70
$ws.Range("H126").Value = 3453.2778
$ws.Range("I126").Value = 3368.7646
$ws.Range("J126").Value = 4890
$ws.Range("K126").Value = 10106.2938
$ws.Range("L126").Value = 14670
$ws.Range("M126").Value = -7636.293799999999
$ws.Range("N126").Value = -19610
$ws.Range("H132").Value = 142858720
$ws.Range("I132").Value = 1731.091
$ws.Range("J132").Value = 666667650
$ws.Range("K132").Value = 5193.272999999999
$ws.Range("L132").Value = 2000002950
$ws.Range("M132").Value = -2663.272999999999
$ws.Range("N132").Value = -2000008010

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3886.3103
$ws.Range("I126").Value = 5577.353
$ws.Range("J126").Value = 1490.6666
$ws.Range("K126").Value = 16732.059
$ws.Range("L126").Value = 4471.9998
$ws.Range("M126").Value = -14262.059
$ws.Range("H136").Value = 2148.5483
$ws.Range("I136").Value = 2201.8462
$ws.Range("J136").Value = 1871.4
$ws.Range("K136").Value = 6605.5386
$ws.Range("L136").Value = 5614.200000000001
$ws.Range("M136").Value = -4055.5386
$ws.Range("N136").Value = -10714.2
